# Update countries & provincias Spain
# - Chequia overtakes Serbia in the ranking (rows 69/70 swap identity; Chequia gets
#   refreshed totals, Serbia keeps its previous totals).
# - Montserrat overtakes Islas Malvinas in the ranking (rows 214/215 swap identity).
# - Refresh case totals for several countries.
# - Bump the "last updated" timestamp caption.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Chequia / Serbia swap places (rows 69-70) ---
$ws.Range("A69").Value = "Chequia"
$ws.Range("B69").Value = 32413
$ws.Range("C69").Value = 1377
$ws.Range("D69").Value = 20787
$ws.Range("E69").Value = 11178
$ws.Range("F69").Value = 0
$ws.Range("G69").Value = 4
$ws.Range("H69").Value = 448

$ws.Range("A70").Value = "Serbia"
$ws.Range("B70").Value = 32136
$ws.Range("C70").Value = 58
$ws.Range("D70").Value = 30943
$ws.Range("E70").Value = 464
$ws.Range("F70").Value = 0
$ws.Range("G70").Value = 1
$ws.Range("H70").Value = 729

# --- Montserrat / Islas Malvinas swap places (rows 214-215) ---
$ws.Range("A214").Value = "Montserrat"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 1

$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 13
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 0

# --- Refresh statistics for other countries ---
# Estados Unidos (row 4)
$ws.Range("B4").Value = 6586895
$ws.Range("C4").Value = 37543
$ws.Range("D4").Value = 3876168
$ws.Range("E4").Value = 2514544
$ws.Range("G4").Value = 945
$ws.Range("H4").Value = 196183

# Brasil (row 6)
$ws.Range("B6").Value = 4239763
$ws.Range("C6").Value = 40431
$ws.Range("E6").Value = 612851
$ws.Range("G6").Value = 922
$ws.Range("H6").Value = 129575

# Peru (row 8)
$ws.Range("B8").Value = 710067
$ws.Range("C8").Value = 7291
$ws.Range("D8").Value = 544745
$ws.Range("E8").Value = 134978
$ws.Range("G8").Value = 108
$ws.Range("H8").Value = 30344

# Argentina (row 13)
$ws.Range("B13").Value = 524198
$ws.Range("C13").Value = 11905
$ws.Range("E13").Value = 123193
$ws.Range("G13").Value = 249
$ws.Range("H13").Value = 10907

# Canada (row 29)
$ws.Range("B29").Value = 134924
$ws.Range("C29").Value = 630
$ws.Range("D29").Value = 118990
$ws.Range("E29").Value = 6771

# Surinam (row 123)
$ws.Range("B123").Value = 4477
$ws.Range("C123").Value = 30
$ws.Range("D123").Value = 3706
$ws.Range("E123").Value = 678
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 93

# Polinesia Francesa (row 169)
$ws.Range("B169").Value = 857
$ws.Range("D169").Value = 576
$ws.Range("E169").Value = 281

# Bermudas (row 190)
$ws.Range("D190").Value = 160
$ws.Range("E190").Value = 8

# --- Bump the "last updated" timestamp caption ---
$ws.Range("A1").Value = "Datos actualizados a 11 de Septiembre de 2020 a las 02:24"
